# Update "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# per-language report sheets, as produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-24 07:56:05"
$zhcn.Range("E4").Value = "2016-03-24 07:56:05"
$zhcn.Range("H2").Value = "2016-03-24 07:56:32"
$zhcn.Range("H4").Value = "2016-03-24 07:56:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-24 07:56:10"
$dede.Range("E4").Value = "2016-03-24 07:56:10"
$dede.Range("H2").Value = "2016-03-24 07:56:39"
$dede.Range("H4").Value = "2016-03-24 07:56:39"
